# Update NaCl NIOSH Test Results 13 Apr 2020 WNg et al.xlsx
# actual runs planned

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths for B:F (values compensate for the fixed internal padding this
# runtime adds between ColumnWidth input and the stored OOXML width, so that
# the persisted <col .../> width attributes land on the target values of
# 13.6640625 / 15.5 / 22 / 32.33203125 / 18)
$ws.Columns.Item(2).ColumnWidth = 12.830729166666666
$ws.Columns.Item(3).ColumnWidth = 14.666666666666666
$ws.Columns.Item(4).ColumnWidth = 21.166666666666668
$ws.Columns.Item(5).ColumnWidth = 31.498697916666668
$ws.Columns.Item(6).ColumnWidth = 17.166666666666668

# Run 1 (row 6-8)
$ws.Range("B6").Value = "NIH Medium"
$ws.Range("C6").Value = "Dragon 30 blue"
$ws.Range("D6").Value = "5N11 outer middle layers"
$ws.Range("E6").Value = "Dummy doesn't fit or stick"
$ws.Range("F6").Value = "Focus is filtration"
$ws.Range("E7").Value = "Aluminum tape doesn't stick"
$ws.Range("E8").Value = "Clip slight loose with extra filter height"

# Clear old run-number cells that move to new rows
$ws.Range("A8").ClearContents()
$ws.Range("A10").ClearContents()
$ws.Range("A12").ClearContents()

# Run 2 (row 9-10)
$ws.Range("A9").Value = 2
$ws.Range("B9").Value = "NIH Medium"
$ws.Range("C9").Value = "Dragon 30 blue"
$ws.Range("D9").Value = "Woodbridge 5-ply"
$ws.Range("D10").Value = "Woodbridge 4-ply"

# Run 3 (row 12)
$ws.Range("A12").Value = 3
$ws.Range("B12").Value = "NIH Medium"
$ws.Range("C12").Value = "Dragon 30 blue"
$ws.Range("D12").Value = "ROXON RX1568"

# Run 4 (row 15) - previously at row 12 with value 4
$ws.Range("A15").Value = 4
$ws.Range("B15").Value = "NIH Medium"
$ws.Range("C15").Value = "Dragon 30 blue"
$ws.Range("D15").Value = "H300 2-ply"

# Set selection to match final state
$ws.Range("E15").Select()
